$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 272, shifting existing rows 272-326 down to 273-327
$ws.Rows.Item(272).Insert()

# Populate the newly inserted row 272 with the new record
$ws.Cells.Item(272, 1).Value = 4
$ws.Cells.Item(272, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(272, 3).Value = "Los Lagos"
$ws.Cells.Item(272, 4).Value = 44782
$ws.Cells.Item(272, 5).Value = 10
$ws.Cells.Item(272, 6).Value = 100112040
$ws.Cells.Item(272, 7).Value = "Cilantro"
$ws.Cells.Item(272, 8).Value = "Sin especificar"
$ws.Cells.Item(272, 9).Value = "Primera"
$ws.Cells.Item(272, 10).Value = 240
$ws.Cells.Item(272, 11).Value = 12500
$ws.Cells.Item(272, 12).Value = 14000
$ws.Cells.Item(272, 13).Value = 13250
$ws.Cells.Item(272, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(272, 15).Value = "Región Metropolitana"
$ws.Cells.Item(272, 16).Value = 368
$ws.Cells.Item(272, 17).Value = 36
$ws.Cells.Item(272, 18).Value = "Hortaliza"
